# Update the South Sudan registration-data-import template to add a
# passport/phone-number column ("pp_phone_no_i_c") with a sample value,
# mirroring the "accomodate phone number in south sudan update script"
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column V: header ------------------------------------------------
$ws.Range("V1").Value = "pp_phone_no_i_c"
$ws.Range("V1").Font.Name = "Arial"
$ws.Range("V1").Font.Size = 10
$ws.Range("V1").Font.Color = 0
$ws.Range("V1").Font.Bold = $false
$ws.Range("V1").Font.Italic = $false

# --- New column V: sample phone number value on the data row ------------
# A leading "+" makes Excel's value parser treat the literal as a number
# (dropping the sign), so build it as a text formula first and then
# collapse it down to a plain stored value via paste-special, which keeps
# it tagged as text.
$ws.Range("V2").Formula = '="+48603499023"'
$ws.Range("V2").Copy()
$ws.Range("V2").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Row 3 (the second sample/data row) is left without a phone number, same
# as the other optional columns in that row.

# --- Header row now matches the shorter row height used elsewhere -------
$ws.Rows.Item(1).RowHeight = 13.8

# --- Move the active selection onto the newly added cell ----------------
$ws.Range("V3").Select()
